$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates (row 2 results from a re-run of the "Cancelacion Anticipada" job) ---
# A2: usuarioAp value used to process the cancelation
$ws.Range("A2").Value = "SCISNEROSC1"
# E2: usuario (customer account holder) associated with the cancelation
$ws.Range("E2").Value = "SCISNEROSCSUP1"
# H2: timestamp of the latest run
$ws.Range("H2").Value = "3 jul. 2023, 16:33:11"

# --- Column A width was manually resized (AutoFit/bestFit turned off) ---
$ws.Columns("A").ColumnWidth = 13

# --- Selection moved to D7 before the file was saved ---
[void]$ws.Range("D7").Select()
